$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of column P (rows 3-10) into the new column Q
# so the new column inherits identical borders/fonts/number formats.
$ws.Range("P3:P10").Copy() | Out-Null
$ws.Range("Q3:Q10").PasteSpecial(-4122) | Out-Null

# New column header: year 2023
$ws.Range("Q4").Value = 2023

# New column data values (mirrors the 2023 figures added alongside 2022's column P)
$ws.Range("Q6").Value = 1209
$ws.Range("Q7").Value = "-"
$ws.Range("Q8").Value = 373
$ws.Range("Q9").Value = 115
$ws.Range("Q10").Value = 781

# Row 5 grew slightly taller to fit the extra column's wrapped header text
$ws.Rows.Item(5).RowHeight = 27
